# Fix weight of tantos
#
# Adds a new "Dagger" flag column (H) to the Stats sheet: TRUE for the
# Dagger/Tanto rows (light, dagger-type weapons) and FALSE for every other
# weapon category, then leaves the Stats sheet as the active tab/selection
# (it had been left on the Patches sheet).

$wb = $excel.ActiveWorkbook
$stats = $wb.Worksheets.Item("Stats")

# New header for column H (reuses the existing "Dagger" shared string).
$stats.Range("H1").Value = "Dagger"

# Dagger (row 2) and Tanto (row 3) are the light dagger-class weapons.
$stats.Range("H2").Value = $true
$stats.Range("H3").Value = $true

# All remaining weapon categories are not daggers.
$stats.Range("H5").Value = $false
$stats.Range("H6").Value = $false
$stats.Range("H7").Value = $false
$stats.Range("H8").Value = $false
$stats.Range("H10").Value = $false
$stats.Range("H11").Value = $false
$stats.Range("H13").Value = $false
$stats.Range("H14").Value = $false
$stats.Range("H15").Value = $false
$stats.Range("H17").Value = $false
$stats.Range("H18").Value = $false
$stats.Range("H20").Value = $false
$stats.Range("H22").Value = $false
$stats.Range("H23").Value = $false
$stats.Range("H25").Value = $false
$stats.Range("H26").Value = $false

# Match the printed page setup recorded alongside the new column.
$stats.PageSetup.PaperSize = 9
$stats.PageSetup.Orientation = 1

# The workbook had been left with the Patches sheet active/selected; move
# the active tab + selection back to Stats (where the edit happened).
$stats.Activate()
$stats.Range("J3").Select()
